$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(30).Insert()

$ws.Range("A30").Value = 9.3
$ws.Range("B30").Value = "Conduct competitive analysis and integrate insights into prototype"
$ws.Range("C30").Value = "Design"
$ws.Range("D30").Value = "Huda"
$ws.Range("E30").Value = 45938
$ws.Range("F30").Value = 45941
$ws.Range("G30").Value = 9.2
$ws.Range("H30").Value = "Start-to-Start"

[void]$ws.Range("B2").Select()
